$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 178, shifting existing rows 178-215 down to 179-216
$ws.Rows(178).Insert()

# Populate the newly inserted row 178 with its data
$ws.Range("A178").Value = 4
$ws.Range("B178").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C178").Value = "Los Lagos"
$ws.Range("D178").Value = 44900
$ws.Range("E178").Value = 10
$ws.Range("F178").Value = "Fruta"
$ws.Range("G178").Value = 100103
$ws.Range("H178").Value = "Frutos de hueso (carozo)"
$ws.Range("I178").Value = 100103004
$ws.Range("J178").Value = "Durazno"
$ws.Range("K178").Value = "Florida King"
$ws.Range("L178").Value = "Primera"
$ws.Range("M178").Value = 400
$ws.Range("N178").Value = 23000
$ws.Range("O178").Value = 24000
$ws.Range("P178").Value = 23500
$ws.Range("Q178").Value = "$/caja 14 kilos empedrada"
$ws.Range("R178").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S178").Value = 1679
$ws.Range("T178").Value = 14
